$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: fill column C with "-" placeholder (K-NN imputation rows with no usable RMSE) ---
$ws.Range("C4").Value = "-"
$ws.Range("C5").Value = "-"
$ws.Range("C7").Value = "-"
$ws.Range("C8").Value = "-"
$ws.Range("C9").Value = "-"
$ws.Range("C10").Value = "-"
$ws.Range("C11").Value = "-"
$ws.Range("C12").Value = "-"
$ws.Range("C14").Value = "-"
$ws.Range("C15").Value = "-"
$ws.Range("C16").Value = "-"
$ws.Range("C17").Value = "-"
$ws.Range("C18").Value = "-"
$ws.Range("C20").Value = "-"
$ws.Range("C21").Value = "-"
$ws.Range("C22").Value = "-"
$ws.Range("C24").Value = "-"
$ws.Range("C25").Value = "-"
$ws.Range("C26").Value = "-"
$ws.Range("C28").Value = "-"
$ws.Range("C29").Value = "-"
$ws.Range("C30").Value = "-"
$ws.Range("C32").Value = "-"
$ws.Range("C33").Value = "-"
$ws.Range("C34").Value = "-"
$ws.Range("C36").Value = "-"
$ws.Range("C37").Value = "-"
$ws.Range("C38").Value = "-"

# --- Step 2: fill column B with actual K-NN output values / discard notices (top block) ---
$ws.Range("B4").Value = "27582.549 (5511.916)"
$ws.Range("B5").Value = "All data have benn discarded"
$ws.Range("B7").Value = "28880.901 (6076.030)"
$ws.Range("B8").Value = "27590.100 (5304.058)"
$ws.Range("B9").Value = "29731.784 (5969.208)"
$ws.Range("B10").Value = "All data have benn discarded"
$ws.Range("B11").Value = "All data have benn discarded"
$ws.Range("B12").Value = "All data have benn discarded"

# --- Step 3: fill column B for rows 56-70 (All data discarded) ---
$ws.Range("B56").Value = "All data have benn discarded"
$ws.Range("B57").Value = "All data have benn discarded"
$ws.Range("B58").Value = "All data have benn discarded"
$ws.Range("B60").Value = "All data have benn discarded"
$ws.Range("B61").Value = "All data have benn discarded"
$ws.Range("B62").Value = "All data have benn discarded"
$ws.Range("B64").Value = "All data have benn discarded"
$ws.Range("B65").Value = "All data have benn discarded"
$ws.Range("B66").Value = "All data have benn discarded"
$ws.Range("B68").Value = "All data have benn discarded"
$ws.Range("B69").Value = "All data have benn discarded"
$ws.Range("B70").Value = "All data have benn discarded"

# --- Step 4: widen column B to fit the new, longer text ---
# (target stored width is 32.33203125 chars; this engine stores ColumnWidth
#  pre-padding and then quantizes to whole screen pixels, so we back out the
#  5/6-character padding before assigning to land as close as possible once
#  Excel re-applies its own padding/quantization on save)
$ws.Range("B1").ColumnWidth = 31.498697916666668

# --- Step 5: move the active selection to C26, matching the author's last edit location ---
$ws.Range("C26").Select()
